# Update the "Förändrad" date column (C2:C7) from 2023-09-05 (45174) to
# 2023-09-06 (45175), leaving everything else untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

foreach ($row in 2..7) {
    $cell = $ws.Cells.Item($row, 3)   # Column C
    if ($cell.Value2 -eq 45174) {
        $cell.Value = 45175
    }
}
